# Inserts a new weekly price record for "Ají" (Terminal Hortofrutícola Agro
# Chillán) at row 52, pushing the existing rows 52-69 down to 53-70. The new
# row gets its own date / price data while the rest of its fields match the
# surrounding records (same market, region, variety, quality, unit, origin,
# classification, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 52..69 down to 53..70, leaving a blank row 52 to fill in.
$ws.Rows.Item(52).Insert()

$ws.Cells.Item(52, 1).Value = 7
$ws.Cells.Item(52, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(52, 3).Value = "Ñuble"
$ws.Cells.Item(52, 4).Value = 44627
$ws.Cells.Item(52, 5).Value = 16
$ws.Cells.Item(52, 6).Value = 100112021
$ws.Cells.Item(52, 7).Value = "Ají"
$ws.Cells.Item(52, 8).Value = "Americana (o)"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 60
$ws.Cells.Item(52, 11).Value = 8500
$ws.Cells.Item(52, 12).Value = 9000
$ws.Cells.Item(52, 13).Value = 8750
$ws.Cells.Item(52, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(52, 15).Value = "Región del Maule"
$ws.Cells.Item(52, 16).Value = 583
$ws.Cells.Item(52, 17).Value = 15
$ws.Cells.Item(52, 18).Value = "Hortaliza"
